# Set start year ("year_start", column B) of construction archetypes to 0
# (was 1950) on the ARCHITECTURE and HVAC sheets, so that buildings with
# unknown age of construction (year = 1) no longer cause the properties
# script to fail. The end year ("year_end", column C, 2030) is unchanged.

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("ARCHITECTURE", "HVAC")) {
    $ws = $wb.Worksheets.Item($sheetName)

    for ($r = 2; $r -le 19; $r++) {
        $ws.Cells.Item($r, 2).Value = "0"
        $ws.Cells.Item($r, 3).Value = "2030"
    }

    # Narrow the remembered selection on the last row to the single cell
    # B19 (it previously spanned B19:C19).
    $ws.Range("B19").Select()
}

$wsArch = $wb.Worksheets.Item("ARCHITECTURE")
$wsArch.Activate()
$wsArch.Range("B19").Select()
